# MVC UPDDATE DetAIL POST
# Update "KEGIATAN BKK" detail-gambar sheet: refresh id_kegiatan / id_kegiatan_gambar
# numbering, and replace/extend the google-drive image hyperlinks for rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the id_kegiatan_gambar (A) and id_kegiatan (B) numeric columns ---
$idKegiatanGambar = @(1036, 1037, 1038, 1039, 1040, 1041, 1042, 1043, 1044, 1045, 1046, 1047, 1048, 1049, 1050)
$idKegiatan = 3617

$startRow = 2
for ($i = 0; $i -lt $idKegiatanGambar.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $idKegiatanGambar[$i]
    $ws.Cells.Item($row, 2).Value = $idKegiatan
}

# --- 2. Replace the "gambar" (C) hyperlink column: drop the old links & rebuild ---
$ws.Hyperlinks.Delete()

$urls = @(
  "http://drive.google.com/uc?export=view&id=1OgocY18GorPFmFv0eLDahUkUYMlUfpT1",
  "http://drive.google.com/uc?export=view&id=1zUNOQxmQaiRCGZ6lcffwJXXFKqZk57gX",
  "http://drive.google.com/uc?export=view&id=1FwvvT4rbtmAPczGsuTpWqIOxmut2x2Z8",
  "http://drive.google.com/uc?export=view&id=18F6ZrsdrhkrmS5qdedS1-was6JY_rIK7",
  "http://drive.google.com/uc?export=view&id=1Gt33V9c8PukZh35uzJVIlZc9HzvE-qfZ",
  "http://drive.google.com/uc?export=view&id=1XwOryjrVcBVcyV0VtJHgskZ1MdPwubGB",
  "http://drive.google.com/uc?export=view&id=10J_GJrYmYl7IJ2G0aGCYdngSy94riw14",
  "http://drive.google.com/uc?export=view&id=1oZ6OI4S6aTHrhQXqZBh6lm7DIGoSSIjL",
  "http://drive.google.com/uc?export=view&id=1xij8W_3LBBMjfgmiJajkHniGuP2w_IzG",
  "http://drive.google.com/uc?export=view&id=1xWv8nrnCBfHXBORNuWR2qqoO5PwlP9bC",
  "http://drive.google.com/uc?export=view&id=1aqV8XnfHWUEYkgxqRuTPyYYk9vhkNT16",
  "http://drive.google.com/uc?export=view&id=1skHeJMFGoAITITsgDGFErVMZxRzPVBEw",
  "http://drive.google.com/uc?export=view&id=16zE4Ov7SwpOESg6whjsrtZ6wk0o_O8tx",
  "http://drive.google.com/uc?export=view&id=1PPBkCUQxRAHgDtn8WGjXKPxVgz3sjZPF",
  "http://drive.google.com/uc?export=view&id=1SyUli3s_yhkKsL7XyCoe0tJhG1xIesGG"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $urls[$i]
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}

# --- 3. Move the active selection, matching the saved sheet view state ---
$ws.Range("A13").Select()
